$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'isolette single sensor thermostat temperature alarm '
$ws.Range("B3").Value = 'security failure datain dataout, isolette dual sensor temp alarm, communication aps architecture dataport'
$ws.Range("B4").Value = 'isolette temp monitoring system with alarm, communication architecture for gnss security, aviate control altitude latitude longitude distance speed, tank pressurized valve switch'
$ws.Range("B5").Value = 'smartparking system for vehicle considering space occupancy and capacity, aviation temperature controller sensor, isolette temperature sensor monitor heat with thermostat and alarm, aps communication architecture'
$ws.Range("B6").Value = 'aviation supervisor for altitude gcas tcas latitude longitude, car cruisecontrol for speed velocity and throttle, blocks memory management using bus'
$ws.Range("B7").Value = 'thesis drone rendezvous control system considering currX currY accX accY velX velY'
$ws.Range("B8").Value = 'radiation pressure sensor shutdown, temp controller heater and cooler, delivery drone system using gps navigation radio,  learning system master slave with write read update sync operations'
$ws.Range("B9").Value = 'aaspe security timing requirements'
$ws.Range("B10").Value = 'smart home remote controller with clients sender and reciever, traffic light signal sensor to switch color'
$ws.Range("B11").Value = 'aaspe security system'
$ws.Range("B12").Value = 'iplprojects tasklibrary map view'
$ws.Range("B13").Value = 'scenario producer consumer'
$ws.Range("B14").Value = 'translator device for abstract and missing data, smartparkingsystem driver '
$ws.Range("B15").Value = 'polyorb scenario types implementation'
$ws.Range("B16").Value = 'flight management scenario'
$ws.Range("B17").Value = 'position control system, cmpare comparing learning model customer service for product and seller, vehicle toll collection payment model, mece learning app system that get req send reply'
$ws.Range("B18").Value = 'vehicle toll collection payment administation system, isolette heat control and monitor system through temperature sensor thermostat and alarm, military radar contact system of target threath position and weapon status, smart parking vehicle system, collision detection and avoidance system'
$ws.Range("B19").Value = 'power manager system with sensorin sensorout actuatorout'
$ws.Range("B20").Value = 'resourcebudgets sensor monitoring hardware power for cpu pci and ram, pcb hardware maintenance system, networking hardware model, automotive car hardware integration cellullar or router for internet access'
$ws.Range("B21").Value = 'aaspe security system'
$ws.Range("B22").Value = 'test impl for bus access'
$ws.Range("B23").Value = 'aaspe security system integration for datain and dataout, flying structural deformation inspection using drones'
$ws.Range("B24").Value = 'structural deformation inspection using drones'
$ws.Range("B25").Value = 'ocarina issue test outp inp'
$ws.Range("B26").Value = 'polyorb rma scenario impl'
$ws.Range("B27").Value = 'communication architecture for networking and system management'
$ws.Range("B28").Value = 'aps communication architecture using lidar and radar'
$ws.Range("B29").Value = 'polyorb ping scenario'
$ws.Range("B30").Value = 'ocarina issue test impl'
$ws.Range("B31").Value = 'ocarina'
$ws.Range("B32").Value = 'iplprojects edge map'
$ws.Range("B33").Value = 'polyorb scenario impl, ocarina system'
$ws.Range("B34").Value = 'polyorb sunseeker scenario impl'
$ws.Range("B35").Value = 'ksu isolette single and dual sensor'
$ws.Range("B36").Value = 'polyorb rma impl test, ocarina rma, sunseekercontrolsystem'
$ws.Range("B37").Value = 'ping pong event system node'
$ws.Range("B38").Value = 'integration model for merged models'
$ws.Range("B39").Value = 'smart home sensor system'
$ws.Range("B40").Value = 'paparazzo ariborne system, vending machine system, car collision detection sensor avoidable unavoidable'
$ws.Range("B41").Value = 'self driving car integration system for speed distance brake and obstacle managing, humidifier and dehumidifier managing sensor'
$ws.Range("B42").Value = 'pacemaker for heart rate monitoring model, radio and gps system, coffeemachine system impl, flightcontroller example impl, groundstation example impl'
$ws.Range("B43").Value = 'pulseox forwarding system, train movement authority controller, flowlatencysampledata application system'
$ws.Range("B44").Value = 'pca management system for patient drug infusion, isolette temperature sensor for heat source, gps satellite position observation, ocarina vehicle speed monitor'
$ws.Range("B45").Value = 'adiru event monitoring alarm system, dca resourcebudget app with sensorin and actuatorout, security algorithm system for crypto, mine pump monitoring software'

$ws.Range("A1:B3").Select()
